$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Calendar Year"
$ws.Range("G1").Value = "Total Firearms"
